$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert new row of data at row 5 (pushing existing row 5.. down isn't needed,
# since the target row 13 stays at row 13 - we just fill row 5 values directly)
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 10000
$ws.Cells.Item(5, 3).Value = "Tan Nguyen"
$ws.Cells.Item(5, 4).Value = 40465
$ws.Cells.Item(5, 5).Value = "late for meeting"
$ws.Cells.Item(5, 6).Value = "Waiting"

# copy style (number format) from D4 to D5 so it matches existing date cells
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# update selection to G4 as in the diff
$ws.Range("G4").Select() | Out-Null
